$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.969.46'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.882.85'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.85'
$ws.Range("E5").Value = '  +8.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.36'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.608'
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.717'
$ws.Range("E9").Value = '  -3.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("E10").Value = '  -5.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000328'
$ws.Range("E11").Value = '  -6.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.01'
$ws.Range("E12").Value = '  -2.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.27'
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.486.36'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.56'
$ws.Range("E15").Value = '  +7.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.848.44'
$ws.Range("E16").Value = '  -3.08%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  +6.32%  '
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.942.23'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '423.70'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.40'
$ws.Range("E22").Value = '  -3.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.14'
$ws.Range("E23").Value = '  -4.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.41'
$ws.Range("E24").Value = '  -3.74%  '
$ws.Range("E25").Value = '  +6.80%  '
$ws.Range("E26").Value = '  -8.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.67'
$ws.Range("E27").Value = '  -3.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.16'
$ws.Range("E28").Value = '  -2.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '689.68'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.15'
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.85'
$ws.Range("E32").Value = '  -2.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '67.03'
$ws.Range("E33").Value = '  +9.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.439'
$ws.Range("E34").Value = '  -2.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.96'
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0854'
$ws.Range("E36").Value = '  -3.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.08'
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0480'
$ws.Range("E42").Value = '  -2.66%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.19'
$ws.Range("E43").Value = '  +5.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  -10.45%  '
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("E46").Value = '  -1.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.817.68'
$ws.Range("E47").Value = '  +16.53%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000275'
$ws.Range("E48").Value = '  +14.26%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.96'
$ws.Range("E49").Value = '  +5.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0346'
$ws.Range("E50").Value = '  -9.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.76'
$ws.Range("E51").Value = '  +1.39%  '
